$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for 2022-Q4 above the existing
#    2022-Q3 row, shifting 2022-Q3 / 2021-Q4 down by one row.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Rows(2).Insert()
$ws1.Range("A2:D2").ClearFormats()
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 15
$ws1.Range("D2").Value = 3.34

# Match the "A" column style used by the other data rows (s="2" in the xf
# table: bordered, bold, centred header-ish number cell).
$ws1.Range("A3").Copy()
$ws1.Range("A2").PasteSpecial(-4122)

# Renumber the index column (0-based) now that a row was inserted.
$ws1.Range("A3").Value = 1
$ws1.Range("A4").Value = 2

# ---------------------------------------------------------------------------
# 2) New "2022-Q4" worksheet: clone the existing "2022-Q3" sheet (so it
#    inherits identical column layout / header styling), place it right
#    after "总计", rename it, then overwrite its data with the Q4 figures.
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Copy($null, $ws1)

$newWs = $wb.Worksheets.Item(2)
$newWs.Name = "2022-Q4"

# The Q3 template sheet has one extra data row (16 funds vs. 15) - drop it.
$newWs.Rows(17).Delete()

# Columns B (fund code) and D:G (scale/position numbers) are stored as text
# in this workbook (e.g. fund codes keep leading zeros) - force text format
# before writing so Excel doesn't silently coerce them to numbers.
$newWs.Range("B2:B16").NumberFormat = "@"
$newWs.Range("D2:G16").NumberFormat = "@"

$newWs.Range("B2").Value = "513500"
$newWs.Range("C2").Value = "博时标普500ETF（QDII）"
$newWs.Range("D2").Value = "93.59"
$newWs.Range("E2").Value = "95.36"
$newWs.Range("F2").Value = "1.47"
$newWs.Range("G2").Value = "1.3758"
$newWs.Range("H2").Value = 6
$newWs.Range("B3").Value = "000041"
$newWs.Range("C3").Value = "华夏全球精选股票（QDII）"
$newWs.Range("D3").Value = "18.44"
$newWs.Range("E3").Value = "85.51"
$newWs.Range("F3").Value = "4.52"
$newWs.Range("G3").Value = "0.8335"
$newWs.Range("H3").Value = 2
$newWs.Range("B4").Value = "000043"
$newWs.Range("C4").Value = "嘉实美国成长股票（QDII）人民币"
$newWs.Range("D4").Value = "12.69"
$newWs.Range("E4").Value = "92.23"
$newWs.Range("F4").Value = "1.80"
$newWs.Range("G4").Value = "0.2284"
$newWs.Range("H4").Value = 6
$newWs.Range("B5").Value = "000044"
$newWs.Range("C5").Value = "嘉实美国成长股票（QDII）美元现汇"
$newWs.Range("D5").Value = "12.69"
$newWs.Range("E5").Value = "92.23"
$newWs.Range("F5").Value = "1.80"
$newWs.Range("G5").Value = "0.2284"
$newWs.Range("H5").Value = 6
$newWs.Range("B6").Value = "000369"
$newWs.Range("C6").Value = "广发全球医疗保健（QDII）人民币A"
$newWs.Range("D6").Value = "3.16"
$newWs.Range("E6").Value = "80.87"
$newWs.Range("F6").Value = "6.63"
$newWs.Range("G6").Value = "0.2095"
$newWs.Range("H6").Value = 1
$newWs.Range("B7").Value = "000370"
$newWs.Range("C7").Value = "广发全球医疗保健（QDII）美元A"
$newWs.Range("D7").Value = "3.16"
$newWs.Range("E7").Value = "80.87"
$newWs.Range("F7").Value = "6.63"
$newWs.Range("G7").Value = "0.2095"
$newWs.Range("H7").Value = 1
$newWs.Range("B8").Value = "012860"
$newWs.Range("C8").Value = "易方达标普500指数（QDII-LOF）人民币 C"
$newWs.Range("D8").Value = "4.75"
$newWs.Range("E8").Value = "91.65"
$newWs.Range("F8").Value = "1.41"
$newWs.Range("G8").Value = "0.0670"
$newWs.Range("H8").Value = 6
$newWs.Range("B9").Value = "161125"
$newWs.Range("C9").Value = "易方达标普500指数（QDII-LOF）人民币"
$newWs.Range("D9").Value = "4.75"
$newWs.Range("E9").Value = "91.65"
$newWs.Range("F9").Value = "1.41"
$newWs.Range("G9").Value = "0.0670"
$newWs.Range("H9").Value = 6
$newWs.Range("B10").Value = "003718"
$newWs.Range("C10").Value = "易方达标普500指数（QDII-LOF）美元A"
$newWs.Range("D10").Value = "4.65"
$newWs.Range("E10").Value = "91.65"
$newWs.Range("F10").Value = "1.41"
$newWs.Range("G10").Value = "0.0656"
$newWs.Range("H10").Value = 6
$newWs.Range("B11").Value = "005698"
$newWs.Range("C11").Value = "华夏全球科技先锋混合（QDII）"
$newWs.Range("D11").Value = "0.60"
$newWs.Range("E11").Value = "83.35"
$newWs.Range("F11").Value = "4.85"
$newWs.Range("G11").Value = "0.0291"
$newWs.Range("H11").Value = 7
$newWs.Range("B12").Value = "159612"
$newWs.Range("C12").Value = "国泰标普500ETF（QDII）"
$newWs.Range("D12").Value = "0.86"
$newWs.Range("E12").Value = "94.21"
$newWs.Range("F12").Value = "1.43"
$newWs.Range("G12").Value = "0.0123"
$newWs.Range("H12").Value = 6
$newWs.Range("B13").Value = "016280"
$newWs.Range("C13").Value = "广发全球医疗保健（QDII）人民币C"
$newWs.Range("D13").Value = "0.10"
$newWs.Range("E13").Value = "80.87"
$newWs.Range("F13").Value = "6.63"
$newWs.Range("G13").Value = "0.0066"
$newWs.Range("H13").Value = 1
$newWs.Range("B14").Value = "016281"
$newWs.Range("C14").Value = "广发全球医疗保健（QDII）美元C"
$newWs.Range("D14").Value = "0.10"
$newWs.Range("E14").Value = "80.87"
$newWs.Range("F14").Value = "6.63"
$newWs.Range("G14").Value = "0.0066"
$newWs.Range("H14").Value = 1
$newWs.Range("B15").Value = "159655"
$newWs.Range("C15").Value = "华夏标普500ETF（QDII）"
$newWs.Range("D15").Value = "0.21"
$newWs.Range("E15").Value = "93.70"
$newWs.Range("F15").Value = "1.44"
$newWs.Range("G15").Value = "0.0030"
$newWs.Range("H15").Value = 5
$newWs.Range("B16").Value = "012861"
$newWs.Range("C16").Value = "易方达标普500指数（QDII-LOF）美元 C"
$newWs.Range("D16").Value = "0.10"
$newWs.Range("E16").Value = "91.65"
$newWs.Range("F16").Value = "1.41"
$newWs.Range("G16").Value = "0.0014"
$newWs.Range("H16").Value = 6
